$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 302..328: column A continues the date series (s="2" style,
# same as the preceding rows), columns B, C, D are all 0.
$startRow = 302
$endRow = 328
$startDate = 44376

for ($r = $startRow; $r -le $endRow; $r++) {
    $dateVal = $startDate + ($r - $startRow)

    # Copy the style (number format, font, border, alignment) from the cell
    # directly above so the new date cell matches the existing "s=2" style.
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))

    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
